# Update countries & provincias Spain
# Applies the diff described for paises.xlsx to the currently-open workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Update the "Datos actualizados..." timestamp cell (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 08:22"

# --- 2. Alemania (row 9): update Casos activos (D) / Recuperados (E) ---
$ws.Range("D9").Value = 120400
$ws.Range("E9").Value = 33198

# --- 3. Bulgaria moves above Cuba (new totals bump it up the ranking) ---
# Row 80 becomes Bulgaria with fresh data
$ws.Range("A80").Value = "Bulgaria"
$ws.Range("C80").Value = 38
$ws.Range("D80").Value = 243
$ws.Range("E80").Value = 1133
$ws.Range("F80").Value = 39
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = 61

# Row 81 becomes Cuba (its previous, unchanged data shifts down one row)
$ws.Range("A81").Value = "Cuba"
$ws.Range("B81").Value = 1437
$ws.Range("D81").Value = 575
$ws.Range("E81").Value = 804
$ws.Range("F81").Value = 12
$ws.Range("H81").Value = 58

# Row 82 becomes Republica de Macedonia (its previous data shifts down one row)
$ws.Range("A82").Value = "Republica de Macedonia"
$ws.Range("B82").Value = 1421
$ws.Range("D82").Value = 589
$ws.Range("E82").Value = 761
$ws.Range("F82").Value = 13
$ws.Range("H82").Value = 71

# Row 83 becomes Eslovenia (its previous data shifts down one row)
$ws.Range("A83").Value = "Eslovenia"
$ws.Range("B83").Value = 1408
$ws.Range("D83").Value = 223
$ws.Range("E83").Value = 1099
$ws.Range("F83").Value = 24
$ws.Range("H83").Value = 86

# Row 84 (Eslovaquia) is unchanged.

# --- 4. Consejo Danes para los Refugiados (row 110): update data ---
$ws.Range("B110").Value = 491
$ws.Range("C110").Value = 20
$ws.Range("D110").Value = 59
$ws.Range("E110").Value = 402

# --- 5. Taiwan (row 114): update Casos activos (D) / Recuperados (E) ---
$ws.Range("D114").Value = 311
$ws.Range("E114").Value = 112

# --- 6. El Salvador moves above Kenia (new totals bump it up the ranking) ---
# Row 118 becomes El Salvador with fresh data
$ws.Range("A118").Value = "El Salvador"
$ws.Range("B118").Value = 377
$ws.Range("C118").Value = 32
$ws.Range("D118").Value = 106
$ws.Range("E118").Value = 262
$ws.Range("F118").Value = 3
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 9

# Row 119 becomes Kenia (its previous, unchanged data shifts down one row)
$ws.Range("A119").Value = "Kenia"
$ws.Range("B119").Value = 374
$ws.Range("D119").Value = 124
$ws.Range("E119").Value = 236
$ws.Range("F119").Value = 2
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 14
